# Implements "UOM column" change:
#   - Column E (on Inputs / Balance Sheet / Corkscrew) is cleared out
#     (value + number format + font reset to plain default), making room
#     for a future UOM column.
#   - A brand new column O is appended, mirroring column N's values /
#     formulas / number formats (shifting any "N" column references to
#     "O").
#
# xlPasteFormats constant (Excel.XlPasteType.xlPasteFormats)
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

function Clear-ColumnE {
    param($ws, $rows)
    foreach ($r in $rows) {
        $cell = $ws.Cells.Item($r, 5)   # column E = 5
        $cell.ClearContents()
        $cell.ClearFormats()
    }
}

function Copy-ColumnNFormatToO {
    param($ws, $rows)
    foreach ($r in $rows) {
        $src = $ws.Cells.Item($r, 14)   # column N = 14
        $dst = $ws.Cells.Item($r, 15)   # column O = 15
        $src.Copy()
        $dst.PasteSpecial($xlPasteFormats)
    }
}

# ---------------------------------------------------------------------
# Sheet: Inputs
# ---------------------------------------------------------------------
$wsInputs = $wb.Worksheets.Item("Inputs")
$inputsRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20)

Clear-ColumnE $wsInputs $inputsRows
Copy-ColumnNFormatToO $wsInputs $inputsRows
foreach ($r in $inputsRows) {
    $wsInputs.Cells.Item($r, 15).Value = 0
}

# ---------------------------------------------------------------------
# Sheet: Balance Sheet
# ---------------------------------------------------------------------
$wsBS = $wb.Worksheets.Item("Balance Sheet")
$bsRows = @(2,3,4,5,6,8,9,10,11,12,14,15,16,17,19,20,21,23,24,25,27,29,30,31,32,34,35,36,37,39,40,41,43,44,45,46,48,49,50,52,53,54,56,57)

Clear-ColumnE $wsBS $bsRows
Copy-ColumnNFormatToO $wsBS $bsRows

$bsFormulas = @{
    2  = "=O30"
    3  = "=Inputs!O`$2"
    4  = "=O2+O3"
    5  = "=Inputs!O`$3"
    6  = "=O5+O4"
    8  = "=Inputs!O`$4"
    9  = "=O32"
    10 = "=O8+O9"
    11 = "=Inputs!O`$5"
    12 = "=O11+O10"
    14 = "=Inputs!O`$6"
    15 = "=Corkscrew!O5"
    16 = "=Inputs!O`$7"
    17 = "=O16+O15+O14"
    19 = "=O17"
    20 = "=O8"
    21 = "=O11"
    23 = "=O3"
    24 = "=O5"
    25 = "=O23+O24"
    27 = "=O19+O20+O21-O25"
    29 = "=max(0, O27)"
    30 = "=O29"
    31 = "=min(O27, 0)"
    32 = "=-1*O31"
    34 = "=O6"
    35 = "=O12"
    36 = "=O17"
    37 = "=O34-(O35+O36)"
    39 = "=Inputs!O`$8"
    40 = "=Inputs!O`$9"
    41 = "=O39-O40"
    43 = "=Inputs!O`$10"
    44 = "=Inputs!O`$11"
    45 = "=Inputs!O`$12"
    46 = "=O41-(O43+O44)+O45"
    48 = "=Inputs!O`$13"
    49 = "=Inputs!O`$14"
    50 = "=O46-(O48+O49)"
    52 = "=Inputs!O`$15"
    53 = "=Inputs!O`$16"
    54 = "=O50-O52+O53"
    56 = "=Inputs!O`$17"
    57 = "=O54-O56"
}
foreach ($r in $bsRows) {
    $wsBS.Cells.Item($r, 15).Formula = $bsFormulas[$r]
}

# ---------------------------------------------------------------------
# Sheet: Corkscrew
# ---------------------------------------------------------------------
$wsCS = $wb.Worksheets.Item("Corkscrew")
$csRows = @(2,3,4,5)

Clear-ColumnE $wsCS $csRows
Copy-ColumnNFormatToO $wsCS $csRows

$csFormulas = @{
    2 = "=Inputs!O`$19"
    3 = "='Balance Sheet'!O57"
    4 = "=Inputs!O`$20"
    5 = "=O2+O3+O4"
}
foreach ($r in $csRows) {
    $wsCS.Cells.Item($r, 15).Formula = $csFormulas[$r]
}
